$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 28

$ws.Cells.Item($row, 1).Value  = "Testmail #17: Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Cells.Item($row, 2).Value  = "Beste [Naam],`nBedankt voor je e-mail. Ik heb de demo op vrijdag om 11:00 uur bij Van Dijk ingepland. Mocht er iets wijzigen of als er nog vragen zijn, laat het me gerust weten.`nMet vriendelijke groet,`n[Jouw Naam]"
$ws.Cells.Item($row, 3).Value  = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Cells.Item($row, 4).Value  = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 5).Value  = "Planning / Afspraak"
$ws.Cells.Item($row, 6).Value  = "2025-07-29 22:04:18"
$ws.Cells.Item($row, 7).Value  = "Ja"
$ws.Cells.Item($row, 8).Value  = "Nee"
$ws.Cells.Item($row, 9).Value  = "Ja"
$ws.Cells.Item($row, 10).Value = "Nee"

$ws.Rows.Item($row).EntireRow.AutoFit()

